$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.351078441845992
$ws.Range("D2").Value = 0.1855931673119104

$ws.Range("C3").Value = -0.2736976606776806
$ws.Range("D3").Value = 0.7859728598146942

$ws.Range("C4").Value = 0.0542385595296184
$ws.Range("D4").Value = 0.9570625570261804

$ws.Range("C5").Value = 0.4746280052461031
$ws.Range("D5").Value = 0.6380866286512168

$ws.Range("C6").Value = -1.898165335894227
$ws.Range("D6").Value = 0.06619028338820687

$ws.Range("C7").Value = -0.9038133190657295
$ws.Range("D7").Value = 0.3724539421796993

$ws.Range("C8").Value = -0.463503805668036
$ws.Range("D8").Value = 0.6459569070323865

$ws.Range("C9").Value = 0.3385275375889127
$ws.Range("D9").Value = 0.7370474099830422

$ws.Range("C10").Value = 0.7339657971996605
$ws.Range("D10").Value = 0.4680011448381567

$ws.Range("C11").Value = 0.616439427927856
$ws.Range("D11").Value = 0.5417102749687377
